$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (Fecha) for rows 52-125: each pair shifts down to the next pair
# (a new most-recent reading is inserted at the top; oldest reading duplicates into new rows 126-127)
$dValues = @{
  52 = 44579
  53 = 44579
  54 = 44349
  55 = 44349
  56 = 44187
  57 = 44187
  58 = 44391
  59 = 44391
  60 = 44433
  61 = 44433
  62 = 44237
  63 = 44237
  64 = 44400
  65 = 44400
  66 = 44350
  67 = 44350
  68 = 44453
  69 = 44453
  70 = 44475
  71 = 44475
  72 = 44523
  73 = 44523
  74 = 44292
  75 = 44292
  76 = 44574
  77 = 44574
  78 = 44168
  79 = 44168
  80 = 44299
  81 = 44299
  82 = 44160
  83 = 44160
  84 = 44308
  85 = 44308
  86 = 44320
  87 = 44320
  88 = 44306
  89 = 44306
  90 = 44316
  91 = 44316
  92 = 44460
  93 = 44460
  94 = 44272
  95 = 44272
  96 = 44313
  97 = 44313
  98 = 44334
  99 = 44334
  100 = 44405
  101 = 44405
  102 = 44280
  103 = 44280
  104 = 44330
  105 = 44330
  106 = 44239
  107 = 44239
  108 = 44476
  109 = 44476
  110 = 44250
  111 = 44250
  112 = 44488
  113 = 44488
  114 = 44341
  115 = 44341
  116 = 44278
  117 = 44278
  118 = 44194
  119 = 44194
  120 = 44490
  121 = 44490
  122 = 44525
  123 = 44525
  124 = 44327
  125 = 44327
}
foreach ($r in $dValues.Keys) {
  $ws.Cells.Item($r, 4).Value = $dValues[$r]
}

# Column O (Origen) follows the same row values as column D for most rows, except the
# following rows where the source region text itself changed between readings
$oValues = @{
  54 = 'Región Metropolitana'
  55 = 'Región Metropolitana'
  108 = 'Región de Ñuble'
  109 = 'Región de Ñuble'
  110 = 'Región de Arica y Parinacota'
  111 = 'Región de Arica y Parinacota'
}
foreach ($r in $oValues.Keys) {
  $ws.Cells.Item($r, 15).Value = $oValues[$r]
}

# Column J (Volumen) follow-up adjustments for the same shifted rows
$jValues = @{
  114 = 200
  115 = 100
  116 = 300
  117 = 150
}
foreach ($r in $jValues.Keys) {
  $ws.Cells.Item($r, 10).Value = $jValues[$r]
}

# New rows 126 and 127: duplicate of the oldest reading (previously rows 124-125)
$ws.Cells.Item(126, 1).Value = 11
$ws.Cells.Item(126, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(126, 3).Value = 'Bíobío'
$ws.Cells.Item(126, 4).Value = 44512
$ws.Cells.Item(126, 5).Value = 8
$ws.Cells.Item(126, 6).Value = 100112044
$ws.Cells.Item(126, 7).Value = 'Perejil'
$ws.Cells.Item(126, 8).Value = 'Sin especificar'
$ws.Cells.Item(126, 9).Value = 'Primera'
$ws.Cells.Item(126, 10).Value = 200
$ws.Cells.Item(126, 11).Value = 600
$ws.Cells.Item(126, 12).Value = 700
$ws.Cells.Item(126, 13).Value = 650
$ws.Cells.Item(126, 14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(126, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(126, 16).Value = 650
$ws.Cells.Item(126, 17).Value = 1
$ws.Cells.Item(126, 18).Value = 'Hortaliza'

$ws.Cells.Item(127, 1).Value = 11
$ws.Cells.Item(127, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(127, 3).Value = 'Bíobío'
$ws.Cells.Item(127, 4).Value = 44512
$ws.Cells.Item(127, 5).Value = 8
$ws.Cells.Item(127, 6).Value = 100112044
$ws.Cells.Item(127, 7).Value = 'Perejil'
$ws.Cells.Item(127, 8).Value = 'Sin especificar'
$ws.Cells.Item(127, 9).Value = 'Segunda'
$ws.Cells.Item(127, 10).Value = 100
$ws.Cells.Item(127, 11).Value = 500
$ws.Cells.Item(127, 12).Value = 500
$ws.Cells.Item(127, 13).Value = 500
$ws.Cells.Item(127, 14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(127, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(127, 16).Value = 500
$ws.Cells.Item(127, 17).Value = 1
$ws.Cells.Item(127, 18).Value = 'Hortaliza'

# Re-apply the datetime number format to the date column for the newly shifted/added rows
$ws.Range("D52:D127").NumberFormat = "YYYY-MM-DD HH:MM:SS"
